$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 25.8370741681734
$ws.Cells.Item(2, 3).Value = 12.27253946530717
$ws.Cells.Item(2, 4).Value = 4.870455300791582
$ws.Cells.Item(2, 6).Value = 49.443679360903
$ws.Cells.Item(2, 7).Value = 3.749477471493676
$ws.Cells.Item(2, 9).Value = 34.23915024873828
$ws.Cells.Item(2, 10).Value = 9.628627763892753
$ws.Cells.Item(2, 12).Value = 12.93162625198395
$ws.Cells.Item(2, 14).Value = 20.64861635271204

$ws.Cells.Item(3, 2).Value = 25.4971301462666
$ws.Cells.Item(3, 3).Value = 11.85683366034749
$ws.Cells.Item(3, 4).Value = 4.852725322498434
$ws.Cells.Item(3, 6).Value = 49.3925844483686
$ws.Cells.Item(3, 7).Value = 3.753552955756449
$ws.Cells.Item(3, 9).Value = 34.28777489731461
$ws.Cells.Item(3, 10).Value = 9.646033925087442
$ws.Cells.Item(3, 12).Value = 12.93218478430718
$ws.Cells.Item(3, 14).Value = 20.72356420005293

$ws.Cells.Item(4, 2).Value = 25.29303564807751
$ws.Cells.Item(4, 3).Value = 11.59801586145939
$ws.Cells.Item(4, 4).Value = 4.841518948194675
$ws.Cells.Item(4, 6).Value = 49.3740028033859
$ws.Cells.Item(4, 7).Value = 3.756185053831625
$ws.Cells.Item(4, 9).Value = 34.32537359873
$ws.Cells.Item(4, 10).Value = 9.657299965328464
$ws.Cells.Item(4, 12).Value = 12.93469308349417
$ws.Cells.Item(4, 14).Value = 20.77160148614177

$ws.Cells.Item(5, 2).Value = 25.2111300004462
$ws.Cells.Item(5, 3).Value = 11.4918326921414
$ws.Cells.Item(5, 4).Value = 4.836871265855961
$ws.Cells.Item(5, 6).Value = 49.36964487037883
$ws.Cells.Item(5, 7).Value = 3.757290403936914
$ws.Cells.Item(5, 9).Value = 34.34263519341332
$ws.Cells.Item(5, 10).Value = 9.662036916769917
$ws.Cells.Item(5, 12).Value = 12.93626032220281
$ws.Cells.Item(5, 14).Value = 20.79168635846224

$ws.Cells.Item(6, 2).Value = 25.19760886541337
$ws.Cells.Item(6, 3).Value = 11.47416369934189
$ws.Cells.Item(6, 4).Value = 4.836094615172012
$ws.Cells.Item(6, 6).Value = 49.36911525832382
$ws.Cells.Item(6, 7).Value = 3.757475928124977
$ws.Cells.Item(6, 9).Value = 34.34561840615994
$ws.Cells.Item(6, 10).Value = 9.662832312357844
$ws.Cells.Item(6, 12).Value = 12.93655349618917
$ws.Cells.Item(6, 14).Value = 20.79505224277045

$ws.Cells.Item(7, 2).Value = 25.29192579023733
$ws.Cells.Item(7, 3).Value = 11.59658646880484
$ws.Cells.Item(7, 4).Value = 4.841456596501977
$ws.Cells.Item(7, 6).Value = 49.37393102142946
$ws.Cells.Item(7, 7).Value = 3.756199828206149
$ws.Cells.Item(7, 9).Value = 34.32559855052714
$ws.Cells.Item(7, 10).Value = 9.657363257965546
$ws.Cells.Item(7, 12).Value = 12.93471201216859
$ws.Cells.Item(7, 14).Value = 20.77187029362409

$ws.Cells.Item(8, 2).Value = 25.71896048556318
$ws.Cells.Item(8, 3).Value = 12.13005262603364
$ws.Cells.Item(8, 4).Value = 4.864408135288499
$ws.Cells.Item(8, 6).Value = 49.42340657310096
$ws.Cells.Item(8, 7).Value = 3.750855848144825
$ws.Cells.Item(8, 9).Value = 34.254304799062
$ws.Cells.Item(8, 10).Value = 9.634509598826945
$ws.Cells.Item(8, 12).Value = 12.93136982248498
$ws.Cells.Item(8, 14).Value = 20.67404029035823

$ws.Cells.Item(9, 2).Value = 26.58865376771442
$ws.Cells.Item(9, 3).Value = 13.14040227071592
$ws.Cells.Item(9, 4).Value = 4.906899830496777
$ws.Cells.Item(9, 6).Value = 49.62190462748784
$ws.Cells.Item(9, 7).Value = 3.741399955217056
$ws.Cells.Item(9, 9).Value = 34.17623580972862
$ws.Cells.Item(9, 10).Value = 9.594263576881042
$ws.Cells.Item(9, 12).Value = 12.94196460329554
$ws.Cells.Item(9, 14).Value = 20.49814151315572

$ws.Cells.Item(10, 2).Value = 27.24122444125764
$ws.Cells.Item(10, 3).Value = 13.85199047889752
$ws.Cells.Item(10, 4).Value = 4.93661856046781
$ws.Cells.Item(10, 6).Value = 49.82943119190393
$ws.Cells.Item(10, 7).Value = 3.735068704765715
$ws.Cells.Item(10, 9).Value = 34.15689441039251
$ws.Cells.Item(10, 10).Value = 9.567451552806101
$ws.Cells.Item(10, 12).Value = 12.96015003839637
$ws.Cells.Item(10, 14).Value = 20.37852315227792

$ws.Cells.Item(11, 2).Value = 27.53976218173822
$ws.Cells.Item(11, 3).Value = 14.16745205571886
$ws.Cells.Item(11, 4).Value = 4.949820037639284
$ws.Cells.Item(11, 6).Value = 49.93714358021076
$ws.Cells.Item(11, 7).Value = 3.732320493543569
$ws.Cells.Item(11, 9).Value = 34.15641825372266
$ws.Cells.Item(11, 10).Value = 9.555846446193007
$ws.Cells.Item(11, 12).Value = 12.97066615654802
$ws.Cells.Item(11, 14).Value = 20.32617061144292

$ws.Cells.Item(12, 2).Value = 27.65294623672998
$ws.Cells.Item(12, 3).Value = 14.28561031749214
$ws.Cells.Item(12, 4).Value = 4.954773955586777
$ws.Cells.Item(12, 6).Value = 49.9798326224605
$ws.Cells.Item(12, 7).Value = 3.731298653145051
$ws.Cells.Item(12, 9).Value = 34.15743900378531
$ws.Cells.Item(12, 10).Value = 9.551536524148529
$ws.Cells.Item(12, 12).Value = 12.97496922201407
$ws.Cells.Item(12, 14).Value = 20.30664097456665

$ws.Cells.Item(13, 2).Value = 27.62856573673742
$ws.Cells.Item(13, 3).Value = 14.26022233624373
$ws.Cells.Item(13, 4).Value = 4.953709042940595
$ws.Cells.Item(13, 6).Value = 49.97055447321711
$ws.Cells.Item(13, 7).Value = 3.731517888421077
$ws.Cells.Item(13, 9).Value = 34.15716568852991
$ws.Cells.Item(13, 10).Value = 9.552460983509116
$ws.Cells.Item(13, 12).Value = 12.97402824093429
$ws.Cells.Item(13, 14).Value = 20.3108339279496

$ws.Cells.Item(14, 2).Value = 27.54907181717289
$ws.Cells.Item(14, 3).Value = 14.17719967552578
$ws.Cells.Item(14, 4).Value = 4.950228502400373
$ws.Cells.Item(14, 6).Value = 49.94061762402787
$ws.Cells.Item(14, 7).Value = 3.732236049075175
$ws.Cells.Item(14, 9).Value = 34.15647813829333
$ws.Cells.Item(14, 10).Value = 9.555490171610293
$ws.Cells.Item(14, 12).Value = 12.97101375313749
$ws.Cells.Item(14, 14).Value = 20.32455798976041

$ws.Cells.Item(15, 2).Value = 27.5003938386936
$ws.Cells.Item(15, 3).Value = 14.12617339369292
$ws.Cells.Item(15, 4).Value = 4.948090699735174
$ws.Cells.Item(15, 6).Value = 49.92252753775794
$ws.Cells.Item(15, 7).Value = 3.732678394151789
$ws.Cells.Item(15, 9).Value = 34.15621352303305
$ws.Cells.Item(15, 10).Value = 9.55735665212555
$ws.Cells.Item(15, 12).Value = 12.96920901888982
$ws.Cells.Item(15, 14).Value = 20.3330027664472

$ws.Cells.Item(16, 2).Value = 27.22173888350734
$ws.Cells.Item(16, 3).Value = 13.8311983272134
$ws.Cells.Item(16, 4).Value = 4.935749461674981
$ws.Cells.Item(16, 6).Value = 49.82265867333064
$ws.Cells.Item(16, 7).Value = 3.735250950821063
$ws.Cells.Item(16, 9).Value = 34.15709338089557
$ws.Cells.Item(16, 10).Value = 9.568221846498183
$ws.Cells.Item(16, 12).Value = 12.95950778354902
$ws.Cells.Item(16, 14).Value = 20.38198589198734

$ws.Cells.Item(17, 2).Value = 27.05114639473814
$ws.Cells.Item(17, 3).Value = 13.64804340750651
$ws.Cells.Item(17, 4).Value = 4.92809748669421
$ws.Cells.Item(17, 6).Value = 49.76479240552693
$ws.Cells.Item(17, 7).Value = 3.736862829566694
$ws.Cells.Item(17, 9).Value = 34.15976799965915
$ws.Cells.Item(17, 10).Value = 9.575038568292769
$ws.Cells.Item(17, 12).Value = 12.9541298081605
$ws.Cells.Item(17, 14).Value = 20.41256266435227

$ws.Cells.Item(18, 2).Value = 26.95319054912997
$ws.Cells.Item(18, 3).Value = 13.5419276165711
$ws.Cells.Item(18, 4).Value = 4.923666338116624
$ws.Cells.Item(18, 6).Value = 49.73276277016472
$ws.Cells.Item(18, 7).Value = 3.737802362415086
$ws.Cells.Item(18, 9).Value = 34.16208949577123
$ws.Cells.Item(18, 10).Value = 9.579015096143927
$ws.Cells.Item(18, 12).Value = 12.95124782507091
$ws.Cells.Item(18, 14).Value = 20.43034382830232

$ws.Cells.Item(19, 2).Value = 26.92005593003465
$ws.Cells.Item(19, 3).Value = 13.50587011387212
$ws.Cells.Item(19, 4).Value = 4.922160870967883
$ws.Cells.Item(19, 6).Value = 49.72213372673033
$ws.Cells.Item(19, 7).Value = 3.738122609400932
$ws.Cells.Item(19, 9).Value = 34.1630098735132
$ws.Cells.Item(19, 10).Value = 9.580371064712081
$ws.Cells.Item(19, 12).Value = 12.95030837399059
$ws.Cells.Item(19, 14).Value = 20.43639762537651

$ws.Cells.Item(20, 2).Value = 27.0692900214975
$ws.Cells.Item(20, 3).Value = 13.66762111934352
$ws.Cells.Item(20, 4).Value = 4.928915147071829
$ws.Cells.Item(20, 6).Value = 49.77082271984699
$ws.Cells.Item(20, 7).Value = 3.736689957481218
$ws.Cells.Item(20, 9).Value = 34.15940219777548
$ws.Cells.Item(20, 10).Value = 9.574307151874104
$ws.Cells.Item(20, 12).Value = 12.95468044763317
$ws.Cells.Item(20, 14).Value = 20.40928762491847

$ws.Cells.Item(21, 2).Value = 27.57241830751227
$ws.Cells.Item(21, 3).Value = 14.20162157988593
$ws.Cells.Item(21, 4).Value = 4.951252044179109
$ws.Cells.Item(21, 6).Value = 49.94935933401229
$ws.Cells.Item(21, 7).Value = 3.732024597323209
$ws.Cells.Item(21, 9).Value = 34.15664746125025
$ws.Cells.Item(21, 10).Value = 9.554598130800024
$ws.Cells.Item(21, 12).Value = 12.97189048895962
$ws.Cells.Item(21, 14).Value = 20.32051890170766

$ws.Cells.Item(22, 2).Value = 27.90197644489286
$ws.Cells.Item(22, 3).Value = 14.54299910140734
$ws.Cells.Item(22, 4).Value = 4.965587295489589
$ws.Cells.Item(22, 6).Value = 50.07711419446791
$ws.Cells.Item(22, 7).Value = 3.729085317105833
$ws.Cells.Item(22, 9).Value = 34.16184926211404
$ws.Cells.Item(22, 10).Value = 9.542210531804157
$ws.Cells.Item(22, 12).Value = 12.98500737069257
$ws.Cells.Item(22, 14).Value = 20.26422301710499

$ws.Cells.Item(23, 2).Value = 27.72605267000593
$ws.Cells.Item(23, 3).Value = 14.36153152903922
$ws.Cells.Item(23, 4).Value = 4.957960222314188
$ws.Cells.Item(23, 6).Value = 50.00792098839337
$ws.Cells.Item(23, 7).Value = 3.730644059028668
$ws.Cells.Item(23, 9).Value = 34.15843102667399
$ws.Cells.Item(23, 10).Value = 9.548777020576281
$ws.Cells.Item(23, 12).Value = 12.97783625984381
$ws.Cells.Item(23, 14).Value = 20.29411231933364

$ws.Cells.Item(24, 2).Value = 27.06108690575204
$ws.Cells.Item(24, 3).Value = 13.65877257357498
$ws.Cells.Item(24, 4).Value = 4.928545582295191
$ws.Cells.Item(24, 6).Value = 49.76809255629277
$ws.Cells.Item(24, 7).Value = 3.736768072941288
$ws.Cells.Item(24, 9).Value = 34.15956513571209
$ws.Cells.Item(24, 10).Value = 9.574637646093853
$ws.Cells.Item(24, 12).Value = 12.95443084964622
$ws.Cells.Item(24, 14).Value = 20.41076764042217

$ws.Cells.Item(25, 2).Value = 26.35054871023435
$ws.Cells.Item(25, 3).Value = 12.87188290227794
$ws.Cells.Item(25, 4).Value = 4.895670555125698
$ws.Cells.Item(25, 6).Value = 49.55734900826845
$ws.Cells.Item(25, 7).Value = 3.743849270564549
$ws.Cells.Item(25, 9).Value = 34.19070439166227
$ws.Cells.Item(25, 10).Value = 9.604664981599749
$ws.Cells.Item(25, 12).Value = 12.93726715132497
$ws.Cells.Item(25, 14).Value = 20.54403075498496

